$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.77959999999999
$ws.Range("E8").Value = 13.35769999999999
$ws.Range("A12").Value = -22.77180000000002
$ws.Range("E12").Value = 12.83549999999999
$ws.Range("E14").Value = 13.64890000000001
$ws.Range("E22").Value = 11.8557
